$d = $word.ActiveDocument

# Paragraph 1: date in the header line changes from 16.06.24 -> 15.06.24
$p0 = $d.Paragraphs.Item(1)
$p0.Range.Text = @'
⚡️🚀המאמר היומי של מייק 15.06.24:⚡️🚀
'@

# Paragraph 2: paper title
$p1 = $d.Paragraphs.Item(2)
$p1.Range.Text = @'
MEDUSA: Simple LLM Inference Acceleration Framework with Multiple Decoding Heads
'@

# Paragraph 3: intro paragraph
$p2 = $d.Paragraphs.Item(3)
$p2.Range.Text = @'
ב 3 הסקירות האחרונות ראינו כמה שיטות איטרטיביות מקבילות, מבוססות על שיטות יאקובי ו- Gauss-Seidel, המנסות להאיץ את מהירות גנרוט הטקסט (decoding) של מודלי שפה. היום נסקור קצרות מאמר המציע גישה אחרת לאותה הבעיה, שגם מבצעת גנרוט מקבילי של טקסט אבל בשיטה 'טיפה' אחרת.
'@

# Paragraph 4: describes the "heads" mechanism
$p3 = $d.Paragraphs.Item(4)
$p3.Range.Text = @'
בגדול המאמר מציע להוסיף ולאמן כמה ״ראשים״ (שכבה לינארית עם סופטמקס) למודל שפה מאומן. מטרתה של כל ראש כזה היא לחזות טקסט לא החל מהטוקן הבא אלא להתחיל לחזות מהטוקן ה-k אחרי הפרומפט (או הטוקן האחרון שנחזה). כלומר בהינתן פרומפט באורך 10 טוקנים הראש מסדר 3 מגנרט טוקנים החל מהטוקן ה-14 בזמן שמודל שפה רגיל חוזה(מגנרט) החל מהטוקן ה-11. הראשים האלו מחוברים לשכבה האחרונה (לפני שכבת החיזוי) של מודל שפה. כלומר הם מפעילים טרנספורמציה לינארית על ייצוג(תלוי קונטקסט) הטוקן המופק על ידי מודל שפה.
'@

# Paragraph 5: training approach paragraph
$p4 = $d.Paragraphs.Item(5)
$p4.Range.Text = @'
המחברים מציעים שתי דרכים לאמן מודל שפה עם הראשים האלו. הדרך הראשונה היא לאמן רק את הראשים כאשר מודל השפה עצמו נותר מוקפא. הדרך השנייה היא לעשות פיין טיון של מודל שפה מאומן (עם LoRa כמובן). במקרה השני הם משלבים את הלוס הסטנדרטי של מודלי שפה עם זה של הראשים האחרים.
'@

# Paragraph 6: inference-time description
$p5 = $d.Paragraphs.Item(6)
$p5.Range.Text = @'
באינפרנס המחברים לוקחים את החיזויים מהראשים השונים (כמה טוקנים החל מטוקן k לכל ראש) של הראשים השונים ומשלבים אותם בצורה דומה ל- beam search (כאן זה קצת יותר מורכב ונקרא tree-search) כדי לקבל את כמה סדרות של טוקנים (המועמדות) שמהן נבנה החיזוי הסופי של מודל שפה. כדי לבחור את התת-סדרות של טוקנים ״הטובות ביותר״ ביותר הם עושים משהו דומה למה שנעשה ב-speculative decoding קלאסי (טיפה יותר מורכב משם ו-rejection sampling בעניין).
'@

# Paragraph 7: "what's the gain" summary sentence
$p6 = $d.Paragraphs.Item(7)
$p6.Range.Text = @'
אז מה הרווח כאן אתם שואלים? שהראשים מופעלים באופן מקבילי ולפעמים בהפעלה אחת שלהם אנו חוזים כמה טוקנים ולא אחד כמו בגנרוט אוטורגרסיבי רגיל.
'@

# Remove the 5 now-obsolete bullet/sentence paragraphs (original paragraphs 8-12)
for ($i = 0; $i -lt 5; $i++) {
    $d.Paragraphs.Item(8).Range.Delete()
}

$pLast = $d.Paragraphs.Item(8)
$pLast.Range.Text = @'
https://arxiv.org/pdf/2401.10774 
'@

Write-Output $d.Paragraphs.Count
